$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text looks like a plain decimal number (e.g. "1.015")
# must be forced to Text format first, otherwise Excel auto-converts the
# input into a floating point number (binary rounding) instead of keeping
# the literal digit string that the source diff expects.

$ws.Range("D2").Value = "27.996.98"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.890.49"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.95"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("E6").Value = "  +1.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4709"
$ws.Range("E7").Value = "  -0.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3946"
$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.79"
$ws.Range("E9").Value = "  -2.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08013"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.017"
$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.80"
$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("D13").Value = "1.897.30"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.985"
$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.175"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("E16").Value = "  +1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06788"
$ws.Range("E17").Value = "  +2.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.01"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001052"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.18"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.014"
$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("D22").Value = "28.011.00"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.500"
$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.98"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.361"
$ws.Range("E25").Value = "  +2.02%  "

$ws.Range("D26").Value = "2.118.45"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.42"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.02"
$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.103"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.492"
$ws.Range("E30").Value = "  -2.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.59"
$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09565"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9653"
$ws.Range("E33").Value = "  -2.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.645"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.354"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.363"
$ws.Range("E36").Value = "  -6.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06123"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02248"
$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.217"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.232"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5964"
$ws.Range("E41").Value = "  -1.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1903"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.33"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.270"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5699"
$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.23"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.946"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06863"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.53"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("E51").Value = "  -0.57%  "
